$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.835000000000001
$ws.Range("B4").Value = 6.612
$ws.Range("D6").Value = -7.509
$ws.Range("B7").Value = 7.181
$ws.Range("D7").Value = -7.569
$ws.Range("B8").Value = 6.861
$ws.Range("D8").Value = -7.419000000000001
$ws.Range("A11").Value = -21.642
$ws.Range("A12").Value = -21.36
$ws.Range("B12").Value = 6.695
$ws.Range("B14").Value = 6.566999999999998
$ws.Range("A15").Value = -21.018
$ws.Range("D19").Value = -8.032000000000002
$ws.Range("D21").Value = -7.222
$ws.Range("B22").Value = 6.964
$ws.Range("D24").Value = -7.531999999999999
$ws.Range("D25").Value = -7.860999999999999
